$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize column E (narrower, from ~33.16 chars down to ~19.33 chars)
$ws.Columns("E").ColumnWidth = 18.5

# Select column F (about to add a new "slides" column for week 1)
$ws.Columns("F").Select()
